$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = -0.01515604076933738
$ws.Range("D2").Value = 0.9880443009350772

# Row 3
$ws.Range("C3").Value = 0.1420142856458639
$ws.Range("D3").Value = 0.8883612359869568

# Row 4
$ws.Range("C4").Value = 2.373651104558985
$ws.Range("D4").Value = 0.02676738831129644

# Row 5
$ws.Range("C5").Value = 1.707003586460345
$ws.Range("D5").Value = 0.1018984378686507

# Row 6
$ws.Range("C6").Value = 0.1401948696939418
$ws.Range("D6").Value = 0.8897816682873154

# Row 7
$ws.Range("C7").Value = 2.514113070475597
$ws.Range("D7").Value = 0.01974835086073656

# Row 8
$ws.Range("C8").Value = 1.60616301320847
$ws.Range("D8").Value = 0.1224975083954221

# Row 9
$ws.Range("C9").Value = 2.619077474768141
$ws.Range("D9").Value = 0.015668333345499

# Row 10
$ws.Range("C10").Value = 2.428883915500699
$ws.Range("D10").Value = 0.02376934690438648
$ws.Range("G10").Value = "Sí"

# Row 11
$ws.Range("C11").Value = -1.089994036502469
$ws.Range("D11").Value = 0.2875112329324598
